$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" header column, matching the formatting of the existing
# header cells (bold font, border, centered alignment) by copying the
# format from the neighboring header cell (G1).
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("H1").Value = "Save"

# New data value for row 2 under the "Save" column.
$ws.Range("H2").Value = 1
